$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 ("Summe") already carries the formatting (cell styles) we need to
# reuse for the new "Median" (row 30) and "Standardabweichung" (row 31)
# summary rows, so copy its format across rather than re-deriving every
# style index by hand.
$ws.Range("C28:N28").Copy()
$ws.Range("C30:N30").PasteSpecial(-4122)
$ws.Range("C28:N28").Copy()
$ws.Range("C31:N31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 30: "Median" label + the J-column median formula.
$ws.Range("C30").Value = "Median"
$ws.Range("J30").Formula = "=MEDIAN(J7:J25)"

# Row 31: "Standardabweichung" label + population-stddev formulas for J, L, M, N.
$ws.Range("C31").Value = "Standardabweichung"
$ws.Range("J31").Formula = "=ROUND(STDEV.P(J7:J25),1)"
$ws.Range("L31").Formula = "=ROUND(STDEV.P(L7:L25),1)"
$ws.Range("M31").Formula = "=ROUND(STDEV.P(M7:M25),1)"
$ws.Range("N31").Formula = "=ROUND(STDEV.P(N7:N25),1)"

# Sheet view: scrolled/selected cell moved, and column C widened (OOXML
# target width 24.7109375; the host quantizes ColumnWidth to a 1/6-char
# grid internally, so 23.8333333 is the input that lands closest to it).
$ws.Columns("C").ColumnWidth = 23.8333333
[void]$ws.Range("O31").Select()
$excel.ActiveWindow.ScrollColumn = 3
